# control_asistencia.xlsx - add the "12" (día 12) attendance column (P)
# mirroring the existing day columns (E..O), marking every student present
# ("p") for rows 3-34, and leave the selection where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill column P (day 12) for every student row with the same "p" marker
# used throughout the rest of the attendance grid (columns E..O).
$ws.Range("P3:P34").Value = "p"

# Move the view/selection to where it ended up after the edit
# (scrolled down a bit, active cell on P35).
$excel.Goto($ws.Range("C15"), $false) | Out-Null
$ws.Range("P35").Select() | Out-Null
